# Auto-generated script applying scheduled market-price refresh values
# to the Sheets workbook (updates currentAveragePrice* / Leve cost & profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1275.1333
$ws.Range("I6").Value = 268
$ws.Range("K6").Value = 804
$ws.Range("M6").Value = -692
$ws.Range("H9").Value = 138.85715
$ws.Range("I9").Value = 95.333336
$ws.Range("J9").Value = 171.5
$ws.Range("K9").Value = 95.333336
$ws.Range("L9").Value = 171.5
$ws.Range("M9").Value = 73.666664
$ws.Range("N9").Value = -509.5
$ws.Range("H80").Value = 1347.2273
$ws.Range("I80").Value = 1631.1177
$ws.Range("J80").Value = 382
$ws.Range("K80").Value = 4893.3531
$ws.Range("L80").Value = 1146
$ws.Range("M80").Value = -3895.3531
$ws.Range("N80").Value = -3142
$ws.Range("H83").Value = 1347.2273
$ws.Range("I83").Value = 1631.1177
$ws.Range("J83").Value = 382
$ws.Range("K83").Value = 14680.0593
$ws.Range("L83").Value = 3438
$ws.Range("M83").Value = -9688.059300000001
$ws.Range("N83").Value = -13422
$ws.Range("H125").Value = 1265
$ws.Range("I125").Value = 1100
$ws.Range("K125").Value = 9900
$ws.Range("M125").Value = -7440
$ws.Range("H137").Value = 24959.12
$ws.Range("I137").Value = 947.85187
$ws.Range("K137").Value = 2843.55561
$ws.Range("M137").Value = -293.5556099999999
$ws.Range("H138").Value = 2003.2716
$ws.Range("I138").Value = 1770.2941
$ws.Range("J138").Value = 2399.3333
$ws.Range("K138").Value = 5310.8823
$ws.Range("L138").Value = 7197.999899999999
$ws.Range("M138").Value = -170.8823000000002
$ws.Range("N138").Value = -17477.9999
$ws.Range("H141").Value = 779337.4
$ws.Range("J141").Value = 9285.333000000001
$ws.Range("L141").Value = 27855.999
$ws.Range("N141").Value = -38215.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3425.1707
$ws.Range("I32").Value = 2985.282
$ws.Range("J32").Value = 12003
$ws.Range("K32").Value = 2985.282
$ws.Range("L32").Value = 12003
$ws.Range("M32").Value = -2698.282
$ws.Range("N32").Value = -12577
$ws.Range("H102").Value = 1881.6666
$ws.Range("I102").Value = 1881.6666
$ws.Range("K102").Value = 1881.6666
$ws.Range("M102").Value = -259.6666
$ws.Range("H109").Value = 58388
$ws.Range("J109").Value = 58388
$ws.Range("L109").Value = 58388
$ws.Range("N109").Value = -61162
$ws.Range("H132").Value = 1822.9149
$ws.Range("I132").Value = 1370.8462
$ws.Range("K132").Value = 4112.5386
$ws.Range("M132").Value = -1582.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1719358
$ws.Range("I86").Value = 2506876.5
$ws.Range("J86").Value = 669333.3
$ws.Range("K86").Value = 2506876.5
$ws.Range("L86").Value = 669333.3
$ws.Range("M86").Value = -2505753.5
$ws.Range("N86").Value = -671579.3
$ws.Range("H89").Value = 1719358
$ws.Range("I89").Value = 2506876.5
$ws.Range("J89").Value = 669333.3
$ws.Range("K89").Value = 12534382.5
$ws.Range("L89").Value = 3346666.5
$ws.Range("M89").Value = -12528766.5
$ws.Range("N89").Value = -3357898.5
$ws.Range("H94").Value = 1084.2667
$ws.Range("J94").Value = 291.6
$ws.Range("L94").Value = 291.6
$ws.Range("N94").Value = -1193.6
$ws.Range("H99").Value = 1553.5
$ws.Range("I99").Value = 1553
$ws.Range("J99").Value = 1553.6666
$ws.Range("K99").Value = 1553
$ws.Range("L99").Value = 1553.6666
$ws.Range("M99").Value = -55
$ws.Range("N99").Value = -4549.6666
$ws.Range("H105").Value = 2273.739
$ws.Range("I105").Value = 2244.8
$ws.Range("K105").Value = 2244.8
$ws.Range("M105").Value = -497.8000000000002
$ws.Range("H134").Value = 6362.75
$ws.Range("I134").Value = 7468.1113
$ws.Range("J134").Value = 3046.6667
$ws.Range("K134").Value = 22404.3339
$ws.Range("L134").Value = 9140.000100000001
$ws.Range("M134").Value = -19869.3339
$ws.Range("N134").Value = -14210.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1896
$ws.Range("I31").Value = 1362.92
$ws.Range("K31").Value = 1362.92
$ws.Range("M31").Value = -1067.92
$ws.Range("H34").Value = 1896
$ws.Range("I34").Value = 1362.92
$ws.Range("K34").Value = 1362.92
$ws.Range("M34").Value = -1160.92
$ws.Range("H122").Value = 3824.2778
$ws.Range("J122").Value = 7208
$ws.Range("L122").Value = 21624
$ws.Range("N122").Value = -26524
$ws.Range("H132").Value = 1583.5454
$ws.Range("I132").Value = 1089.1708
$ws.Range("K132").Value = 3267.512400000001
$ws.Range("M132").Value = -737.5124000000005
$ws.Range("H134").Value = 1451.6349
$ws.Range("I134").Value = 1356.3617
$ws.Range("J134").Value = 1731.5
$ws.Range("K134").Value = 4069.0851
$ws.Range("L134").Value = 5194.5
$ws.Range("M134").Value = -1534.0851
$ws.Range("N134").Value = -10264.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 94.92308
$ws.Range("I2").Value = 118.333336
$ws.Range("K2").Value = 710.000016
$ws.Range("M2").Value = -597.000016
$ws.Range("H33").Value = 88.14286
$ws.Range("J33").Value = 154.5
$ws.Range("L33").Value = 927
$ws.Range("N33").Value = -1493
$ws.Range("H56").Value = 7374.067
$ws.Range("I56").Value = 7374.067
$ws.Range("K56").Value = 7374.067
$ws.Range("M56").Value = -6844.067
$ws.Range("H131").Value = 13535823
$ws.Range("J131").Value = 30389.74
$ws.Range("L131").Value = 91169.22
$ws.Range("N131").Value = -101249.22

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3045.6538
$ws.Range("I70").Value = 2895.85
$ws.Range("K70").Value = 2895.85
$ws.Range("M70").Value = -2625.85
$ws.Range("H73").Value = 3045.6538
$ws.Range("I73").Value = 2895.85
$ws.Range("K73").Value = 2895.85
$ws.Range("M73").Value = -1959.85
$ws.Range("H97").Value = 978.9259
$ws.Range("I97").Value = 1084.4375
$ws.Range("J97").Value = 825.4545000000001
$ws.Range("K97").Value = 1084.4375
$ws.Range("L97").Value = 825.4545000000001
$ws.Range("M97").Value = -588.4375
$ws.Range("N97").Value = -1817.4545
$ws.Range("H102").Value = 2267.037
$ws.Range("I102").Value = 2196.2104
$ws.Range("J102").Value = 2435.25
$ws.Range("K102").Value = 2196.2104
$ws.Range("L102").Value = 2435.25
$ws.Range("M102").Value = -574.2103999999999
$ws.Range("N102").Value = -5679.25
$ws.Range("H126").Value = 2264761
$ws.Range("I126").Value = 9262505
$ws.Range("K126").Value = 27787515
$ws.Range("M126").Value = -27785045
$ws.Range("H132").Value = 918098
$ws.Range("I132").Value = 1540407.9
$ws.Range("K132").Value = 4621223.699999999
$ws.Range("M132").Value = -4618693.699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 11477.111
$ws.Range("I16").Value = 11477.111
$ws.Range("K16").Value = 11477.111
$ws.Range("M16").Value = -11307.111
$ws.Range("H68").Value = 1588.6471
$ws.Range("I68").Value = 1238.25
$ws.Range("K68").Value = 1238.25
$ws.Range("M68").Value = -489.25
$ws.Range("H71").Value = 1588.6471
$ws.Range("I71").Value = 1238.25
$ws.Range("K71").Value = 6191.25
$ws.Range("M71").Value = -2447.25
$ws.Range("H82").Value = 2345.3333
$ws.Range("I82").Value = 1525.5
$ws.Range("J82").Value = 3985
$ws.Range("K82").Value = 1525.5
$ws.Range("L82").Value = 3985
$ws.Range("M82").Value = -1164.5
$ws.Range("N82").Value = -4707
$ws.Range("H85").Value = 2345.3333
$ws.Range("I85").Value = 1525.5
$ws.Range("J85").Value = 3985
$ws.Range("K85").Value = 1525.5
$ws.Range("L85").Value = 3985
$ws.Range("M85").Value = -277.5
$ws.Range("N85").Value = -6481
$ws.Range("H93").Value = 1435.8182
$ws.Range("I93").Value = 1000.25
$ws.Range("J93").Value = 1684.7142
$ws.Range("K93").Value = 1000.25
$ws.Range("L93").Value = 1684.7142
$ws.Range("M93").Value = 247.75
$ws.Range("N93").Value = -4180.7142
$ws.Range("H132").Value = 3510.3845
$ws.Range("I132").Value = 2226.5293
$ws.Range("J132").Value = 5935.4443
$ws.Range("K132").Value = 6679.5879
$ws.Range("L132").Value = 17806.3329
$ws.Range("M132").Value = -4149.5879
$ws.Range("N132").Value = -22866.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49484.8
$ws.Range("J46").Value = 49484.8
$ws.Range("L46").Value = 49484.8
$ws.Range("N46").Value = -49946.8
$ws.Range("H96").Value = 10754.444
$ws.Range("I96").Value = 3497
$ws.Range("K96").Value = 3497
$ws.Range("M96").Value = -2124
$ws.Range("H132").Value = 1894.1351
$ws.Range("I132").Value = 1439.6086
$ws.Range("K132").Value = 4318.825800000001
$ws.Range("M132").Value = -1788.825800000001
$ws.Range("H134").Value = 49484.8
$ws.Range("J134").Value = 49484.8
$ws.Range("L134").Value = 148454.4
$ws.Range("N134").Value = -153524.4
$ws.Range("H136").Value = 12079540
$ws.Range("I136").Value = 15434379
$ws.Range("J136").Value = 2121
$ws.Range("K136").Value = 46303137
$ws.Range("L136").Value = 6363
$ws.Range("M136").Value = -46300587
$ws.Range("N136").Value = -11463
